$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the driving value; dependent formulas (C3, D3, D4, C8, D8) will recalc automatically.
$ws.Range("C4").Value = 2.59

# Update the active cell selection shown in the sheet view.
$ws.Range("C5").Select()
